$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to match new content width (raw OOXML width=86;
# the COM ColumnWidth property adds ~0.8333 padding when round-tripped)
$ws.Columns.Item(1).ColumnWidth = 85.16666666666666

$newRows = @(
    ,@("What is the maximum number of tracks that can be specified within a single ODF file?", "llama3.2:latest", "Unfortunately, I couldn't find any information on the maximum number of tracks that can be specified within a single ODF file in the provided documentation.")
    ,@("How many tracks can you set up in one ODF?", "llama3.2:latest", "You can set up a maximum of 200 tracks in one ODF.")
    ,@("What’s the track limit for an ODF?", "llama3.2:latest", "The track limit for an ODF (Original Depth File) in GEO is not explicitly stated, but it can be inferred that there are limitations on the number of modifiers, lithologies, symbols, texts, and lines used in an ODF file. If these limits are exceeded, the system may beep and/or display an error message.`nHowever, I couldn't find any specific information on a `"track limit`" for an ODF in the provided documentation.")
    ,@("In one ODF configuration, how many tracks can be defined?", "llama3.2:latest", "According to the GEO application documentation, there is no specific limit mentioned on the number of tracks that can be defined in an ODF configuration. However, it is recommended to keep the track templates organized and consistent for better management and sharing purposes.")
    ,@("Is there a maximum number of tracks you can create in one ODF?", "llama3.2:latest", "Yes, according to the document, there is a limit on the number of tables that can be put in one ODF file. The exact number is not specified, but it mentions that the software has limits on how many tables can be included in an ODF file.")
)

$startRow = 108
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}
